$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("F6").Value = ""
